$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2799.2222
$ws.Range("I135").Value = 2343.3333
$ws.Range("J135").Value = 3027.1667
$ws.Range("K135").Value = 21089.9997
$ws.Range("L135").Value = 27244.5003
$ws.Range("M135").Value = -18554.9997
$ws.Range("N135").Value = -32314.5003
$ws.Range("H137").Value = 2756.8096
$ws.Range("I137").Value = 2772.2778
$ws.Range("K137").Value = 8316.8334
$ws.Range("M137").Value = -5766.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3402.2
$ws.Range("I45").Value = 2340.6667
$ws.Range("J45").Value = 4994.5
$ws.Range("K45").Value = 2340.6667
$ws.Range("L45").Value = 4994.5
$ws.Range("M45").Value = -1963.6667
$ws.Range("N45").Value = -5748.5
$ws.Range("H110").Value = 2977179.2
$ws.Range("I110").Value = 3969119.8
$ws.Range("J110").Value = 1358
$ws.Range("K110").Value = 3969119.8
$ws.Range("L110").Value = 1358
$ws.Range("M110").Value = -3967074.8
$ws.Range("N110").Value = -5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 173.625
$ws.Range("I5").Value = 164.15384
$ws.Range("J5").Value = 214.66667
$ws.Range("K5").Value = 164.15384
$ws.Range("L5").Value = 214.66667
$ws.Range("M5").Value = -51.15384
$ws.Range("N5").Value = -440.66667
$ws.Range("H26").Value = 10955.5
$ws.Range("I26").Value = 10955.5
$ws.Range("K26").Value = 10955.5
$ws.Range("M26").Value = -10663.5
$ws.Range("H94").Value = 2752.875
$ws.Range("I94").Value = 2744.7144
$ws.Range("J94").Value = 2810
$ws.Range("K94").Value = 2744.7144
$ws.Range("L94").Value = 2810
$ws.Range("M94").Value = -2293.7144
$ws.Range("N94").Value = -3712
$ws.Range("H107").Value = 2779.2
$ws.Range("I107").Value = 974.8333
$ws.Range("J107").Value = 3552.5
$ws.Range("K107").Value = 974.8333
$ws.Range("L107").Value = 3552.5
$ws.Range("M107").Value = 945.1667
$ws.Range("N107").Value = -7392.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 817
$ws.Range("J2").Value = 784.6667
$ws.Range("L2").Value = 784.6667
$ws.Range("N2").Value = -1010.6667
$ws.Range("H3").Value = 2500
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 4500
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = -387
$ws.Range("N3").Value = -4726
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0
$ws.Range("H10").Value = 2600
$ws.Range("I10").Value = 2600
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2600
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -2461
$ws.Range("H14").Value = 1772.75
$ws.Range("J14").Value = 2030.3334
$ws.Range("L14").Value = 2030.3334
$ws.Range("N14").Value = -2370.3334
$ws.Range("H15").Value = 921.6667
$ws.Range("I15").Value = 595
$ws.Range("J15").Value = 1085
$ws.Range("K15").Value = 595
$ws.Range("L15").Value = 1085
$ws.Range("M15").Value = -425
$ws.Range("N15").Value = -1425
$ws.Range("H31").Value = 2964.6
$ws.Range("I31").Value = 1193.5
$ws.Range("K31").Value = 1193.5
$ws.Range("M31").Value = -898.5
$ws.Range("H32").Value = 1114.1428
$ws.Range("I32").Value = 1532.25
$ws.Range("J32").Value = 556.6667
$ws.Range("K32").Value = 1532.25
$ws.Range("L32").Value = 556.6667
$ws.Range("M32").Value = -1216.25
$ws.Range("N32").Value = -1188.6667
$ws.Range("H34").Value = 2964.6
$ws.Range("I34").Value = 1193.5
$ws.Range("K34").Value = 1193.5
$ws.Range("M34").Value = -991.5
$ws.Range("H43").Value = 10769.625
$ws.Range("J43").Value = 11593.857
$ws.Range("L43").Value = 11593.857
$ws.Range("N43").Value = -11961.857
$ws.Range("H48").Value = 18922.5
$ws.Range("J48").Value = 18922.5
$ws.Range("L48").Value = 18922.5
$ws.Range("N48").Value = -19874.5
$ws.Range("H54").Value = 22630.666
$ws.Range("J54").Value = 22630.666
$ws.Range("L54").Value = 22630.666
$ws.Range("N54").Value = -23946.666
$ws.Range("H101").Value = 10769.625
$ws.Range("J101").Value = 11593.857
$ws.Range("L101").Value = 11593.857
$ws.Range("N101").Value = -18083.857
$ws.Range("H103").Value = 35140
$ws.Range("I103").Value = 26186.834
$ws.Range("K103").Value = 26186.834
$ws.Range("M103").Value = -25014.834
$ws.Range("H134").Value = 5160.2856
$ws.Range("I134").Value = 5160.2856
$ws.Range("K134").Value = 15480.8568
$ws.Range("M134").Value = -12945.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 216.76471
$ws.Range("I12").Value = 209.66667
$ws.Range("J12").Value = 220.63637
$ws.Range("K12").Value = 629.00001
$ws.Range("L12").Value = 661.9091100000001
$ws.Range("M12").Value = -456.00001
$ws.Range("N12").Value = -1007.90911
$ws.Range("H107").Value = 786.1177
$ws.Range("J107").Value = 1054.4546
$ws.Range("L107").Value = 3163.3638
$ws.Range("N107").Value = -7003.3638
$ws.Range("H137").Value = 4173
$ws.Range("J137").Value = 3963
$ws.Range("L137").Value = 11889
$ws.Range("N137").Value = -22089

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 2500
$ws.Range("I4").Value = 2500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2500
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = -2388
$ws.Range("N4").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 3500
$ws.Range("J14").Value = 3500
$ws.Range("L14").Value = 3500
$ws.Range("N14").Value = -3844
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H40").Value = 7245.467
$ws.Range("I40").Value = 7626.091
$ws.Range("K40").Value = 7626.091
$ws.Range("M40").Value = -7490.091
$ws.Range("H46").Value = 2672.5
$ws.Range("J46").Value = 4900
$ws.Range("L46").Value = 4900
$ws.Range("N46").Value = -5276
$ws.Range("H55").Value = 636.125
$ws.Range("J55").Value = 466.33334
$ws.Range("L55").Value = 466.33334
$ws.Range("N55").Value = -812.33334
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 5165.8335
$ws.Range("I132").Value = 4998.6665
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 14995.9995
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -12465.9995
$ws.Range("N132").Value = -21059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 16113.6
$ws.Range("I21").Value = 500
$ws.Range("K21").Value = 500
$ws.Range("M21").Value = -265
$ws.Range("H25").Value = 24999
$ws.Range("J25").Value = 24999
$ws.Range("L25").Value = 24999
$ws.Range("N25").Value = -25585
$ws.Range("H35").Value = 16113.6
$ws.Range("I35").Value = 500
$ws.Range("K35").Value = 500
$ws.Range("M35").Value = -210
$ws.Range("H122").Value = 1358.8572
$ws.Range("I122").Value = 1338.8
$ws.Range("K122").Value = 4016.4
$ws.Range("M122").Value = -1566.4
